$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'25.877.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.05%  "

# Row 3
$ws.Range("D3").Value = "'1.637.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'215.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("D6").Value = "'0.5033"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.63%  "

# Row 7
$ws.Range("D7").Value = "'1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.40%  "

# Row 8
$ws.Range("D8").Value = "'0.2562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.14%  "

# Row 9
$ws.Range("D9").Value = "'0.06387"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.64%  "

# Row 10
$ws.Range("D10").Value = "'19.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.41%  "

# Row 11
$ws.Range("E11").Value = "  -0.60%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.258"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.69%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.640.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.83%  "

# Row 14
$ws.Range("D14").Value = "'1.863.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.76%  "

# Row 15
$ws.Range("D15").Value = "'0.5441"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.33%  "

# Row 16
$ws.Range("D16").Value = "'0.0₅7905"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "

# Row 17
$ws.Range("D17").Value = "'64.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.26%  "

# Row 18
$ws.Range("D18").Value = "'25.897.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.01%  "

# Row 19
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "

# Row 20
$ws.Range("D20").Value = "'202.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.38%  "

# Row 21
$ws.Range("D21").Value = "'4.373"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22
$ws.Range("D22").Value = "'9.894"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "

# Row 23
$ws.Range("D23").Value = "'5.981"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.21%  "

# Row 24
$ws.Range("D24").Value = "'1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "

# Row 25
$ws.Range("D25").Value = "'1.925"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.19%  "

# Row 26
$ws.Range("D26").Value = "'141.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.96%  "

# Row 27
$ws.Range("D27").Value = "'0.1137"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.19%  "

# Row 28
$ws.Range("D28").Value = "'15.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.83%  "

# Row 29
$ws.Range("D29").Value = "'6.730"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.36%  "

# Row 30
$ws.Range("D30").Value = "'1.242"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("D31").Value = "'0.04954"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.70%  "

# Row 32
$ws.Range("D32").Value = "'3.272"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.40%  "

# Row 33
$ws.Range("D33").Value = "'3.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.09%  "

# Row 34
$ws.Range("D34").Value = "'1.541"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.16%  "

# Row 35
$ws.Range("D35").Value = "'2.372"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "

# Row 36
$ws.Range("E36").Value = "  -3.95%  "

# Row 37
$ws.Range("D37").Value = "'0.8903"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.60%  "

# Row 38
$ws.Range("D38").Value = "'1.159.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.20%  "

# Row 39
$ws.Range("D39").Value = "'0.5596"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.84%  "

# Row 40
$ws.Range("E40").Value = "  -1.37%  "

# Row 41
$ws.Range("D41").Value = "'1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "

# Row 42
$ws.Range("D42").Value = "'5.662"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "

# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.48%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8064"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.15%  "

# Row 45
$ws.Range("D45").Value = "'1.776.17"
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "  +0.63%  "

# Row 47
$ws.Range("D47").Value = "'0.4534"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.40%  "

# Row 48
$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.51%  "

# Row 49
$ws.Range("D49").Value = "'54.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.04%  "

# Row 50
$ws.Range("D50").Value = "'0.05054"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "

# Row 51
$ws.Range("D51").Value = "'1.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
